# Updated cryptos list on Fri Jun  7 08:31:15 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "71.106.93"
    "E2"  = "  +0.23%  "
    "D3"  = "3.811.07"
    "E3"  = "  -0.95%  "
    "E4"  = "  -0.05%  "
    "D5"  = "698.23"
    "E5"  = "  -0.84%  "
    "D6"  = "171.87"
    "E6"  = "  -0.40%  "
    "D7"  = "3.809.68"
    "E7"  = "  -0.97%  "
    "E8"  = "  +0.00%  "
    "D9"  = "0.524"
    "E9"  = "  +0.03%  "
    "E10" = "  -0.47%  "
    "E11" = "  +2.48%  "
    "D12" = "0.464"
    "E12" = "  +1.17%  "
    "E13" = "  -0.92%  "
    "D14" = "35.99"
    "E14" = "  -1.18%  "
    "D15" = "4.451.97"
    "E15" = "  -0.89%  "
    "D16" = "3.829.91"
    "E16" = "  -0.18%  "
    "D17" = "71.057.95"
    "E17" = "  +0.17%  "
    "D18" = "17.47"
    "E18" = "  +0.58%  "
    "E19" = "  -0.41%  "
    "E20" = "  -0.52%  "
    "D21" = "512.72"
    "E21" = "  +3.94%  "
    "E22" = "  -0.22%  "
    "D23" = "0.715"
    "E23" = "  -0.01%  "
    "D24" = "84.07"
    "E24" = "  -1.09%  "
    "E25" = "  -1.28%  "
    "D26" = "3.963.59"
    "E26" = "  -0.88%  "
    "D27" = "12.06"
    "E27" = "  -0.65%  "
    "E28" = "  -1.12%  "
    "E29" = "  +0.13%  "
    "D30" = "2.04"
    "E30" = "  -3.84%  "
    "E31" = "  -4.82%  "
    "E32" = "  -1.08%  "
    "E33" = "  -0.92%  "
    "D34" = "29.09"
    "E34" = "  -1.14%  "
    "D35" = "0.176"
    "E35" = "  -3.81%  "
    "E36" = "  +0.13%  "
    "D37" = "3.773.28"
    "E37" = "  -0.69%  "
    "D38" = "0.999"
    "E38" = "  -0.17%  "
    "D40" = "2.40"
    "E40" = "  +0.85%  "
    "E41" = "  -0.32%  "
    "E42" = "  -0.87%  "
    "E43" = "  +0.18%  "
    "E44" = "  -0.02%  "
    "D45" = "173.17"
    "E45" = "  +5.58%  "
    "E46" = "  +0.08%  "
    "D47" = "0.000310"
    "E47" = "  -0.36%  "
    "E48" = "  +1.58%  "
    "D49" = "429.06"
    "E49" = "  +4.07%  "
    "D50" = "8.65"
    "E50" = "  +0.31%  "
    "B51" = "SuiNetwork"
    "C51" = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
    "D51" = "1.16"
    "E51" = "  +5.40%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "0.000310", "2.40",
    # "71.106.93") keep their exact original formatting instead of being
    # coerced into doubles/dates by Excel's auto-detect.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
